$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: update recalculated TPM-derived values
$ws.Range("G2").Value = 0.8200746666666667
$ws.Range("H2").Value = 2.460224
$ws.Range("I2").Value = 0.04311293902675128
$ws.Range("J2").Value = 0.04311293902675128
$ws.Range("M2").Value = 4.043133999999999
$ws.Range("N2").Value = 12.129402
$ws.Range("O2").Value = 0.7517044794313785
$ws.Range("P2").Value = 0.7517044794313784
$ws.Range("Q2").Value = 3.315671767338666
$ws.Range("R2").Value = 29.841045906048
$ws.Range("S2").Value = 0.03240818938786084
$ws.Range("T2").Value = 0.03240818938786083

# Row 3: update recalculated TPM-derived values
$ws.Range("G3").Value = 0.8200746666666667
$ws.Range("H3").Value = 2.460224
$ws.Range("I3").Value = 0.04311293902675128
$ws.Range("J3").Value = 0.04311293902675128
$ws.Range("M3").Value = 0.4282866666666667
$ws.Range("N3").Value = 1.28486
$ws.Range("O3").Value = 0.07962758736516451
$ws.Range("P3").Value = 0.07962758736516451
$ws.Range("Q3").Value = 0.3512270454044445
$ws.Range("R3").Value = 3.16104340864
$ws.Range("S3").Value = 0.003432979318921648
$ws.Range("T3").Value = 0.003432979318921648

# Row 4: update recalculated TPM-derived values
$ws.Range("G4").Value = 0.8200746666666667
$ws.Range("H4").Value = 2.460224
$ws.Range("I4").Value = 0.04311293902675128
$ws.Range("J4").Value = 0.04311293902675128
$ws.Range("M4").Value = 0.3522683333333333
$ws.Range("N4").Value = 1.056805
$ws.Range("O4").Value = 0.06549416470700518
$ws.Range("P4").Value = 0.06549416470700517
$ws.Range("Q4").Value = 0.2888863360355556
$ws.Range("R4").Value = 2.59997702432
$ws.Range("S4").Value = 0.002823645929621121
$ws.Range("T4").Value = 0.00282364592962112

# Row 5: update recalculated TPM-derived values
$ws.Range("G5").Value = 0.8200746666666667
$ws.Range("H5").Value = 2.460224
$ws.Range("I5").Value = 0.04311293902675128
$ws.Range("J5").Value = 0.04311293902675128
$ws.Range("M5").Value = 0.5549326666666667
$ws.Range("N5").Value = 1.664798
$ws.Range("O5").Value = 0.1031737684964519
$ws.Range("P5").Value = 0.1031737684964519
$ws.Range("Q5").Value = 0.4550862216391112
$ws.Range("R5").Value = 4.095775994752
$ws.Range("S5").Value = 0.004448124390347682
$ws.Range("T5").Value = 0.004448124390347682

# Row 6: update recalculated TPM-derived values
$ws.Range("G6").Value = 1.358031333333334
$ws.Range("H6").Value = 4.074094000000001
$ws.Range("I6").Value = 0.07139437962203982
$ws.Range("J6").Value = 0.07139437962203982
$ws.Range("M6").Value = 4.043133999999999
$ws.Range("N6").Value = 12.129402
$ws.Range("O6").Value = 0.7517044794313785
$ws.Range("P6").Value = 0.7517044794313784
$ws.Range("Q6").Value = 5.490702656865333
$ws.Range("R6").Value = 49.41632391178801
$ws.Range("S6").Value = 0.05366747496811166
$ws.Range("T6").Value = 0.05366747496811165

# Row 7: update recalculated TPM-derived values
$ws.Range("G7").Value = 1.358031333333334
$ws.Range("H7").Value = 4.074094000000001
$ws.Range("I7").Value = 0.07139437962203982
$ws.Range("J7").Value = 0.07139437962203982
$ws.Range("M7").Value = 0.4282866666666667
$ws.Range("N7").Value = 1.28486
$ws.Range("O7").Value = 0.07962758736516451
$ws.Range("P7").Value = 0.07962758736516451
$ws.Range("Q7").Value = 0.5816267129822223
$ws.Range("R7").Value = 5.234640416840001
$ws.Range("S7").Value = 0.005684962200735696
$ws.Range("T7").Value = 0.005684962200735696

# Row 8: update recalculated TPM-derived values
$ws.Range("G8").Value = 1.358031333333334
$ws.Range("H8").Value = 4.074094000000001
$ws.Range("I8").Value = 0.07139437962203982
$ws.Range("J8").Value = 0.07139437962203982
$ws.Range("M8").Value = 0.3522683333333333
$ws.Range("N8").Value = 1.056805
$ws.Range("O8").Value = 0.06549416470700518
$ws.Range("P8").Value = 0.06549416470700517
$ws.Range("Q8").Value = 0.4783914344077779
$ws.Range("R8").Value = 4.305522909670001
$ws.Range("S8").Value = 0.00467591525812033
$ws.Range("T8").Value = 0.004675915258120329

# Row 9: update recalculated TPM-derived values
$ws.Range("G9").Value = 1.358031333333334
$ws.Range("H9").Value = 4.074094000000001
$ws.Range("I9").Value = 0.07139437962203982
$ws.Range("J9").Value = 0.07139437962203982
$ws.Range("M9").Value = 0.5549326666666667
$ws.Range("N9").Value = 1.664798
$ws.Range("O9").Value = 0.1031737684964519
$ws.Range("P9").Value = 0.1031737684964519
$ws.Range("Q9").Value = 0.7536159492235558
$ws.Range("R9").Value = 6.782543543012002
$ws.Range("S9").Value = 0.007366027195072138
$ws.Range("T9").Value = 0.007366027195072136

# Row 10: update recalculated TPM-derived values
$ws.Range("G10").Value = 16.8273
$ws.Range("H10").Value = 50.4819
$ws.Range("I10").Value = 0.8846442749337277
$ws.Range("J10").Value = 0.8846442749337278
$ws.Range("M10").Value = 4.043133999999999
$ws.Range("N10").Value = 12.129402
$ws.Range("O10").Value = 0.7517044794313785
$ws.Range("P10").Value = 0.7517044794313784
$ws.Range("Q10").Value = 68.03502875819999
$ws.Range("R10").Value = 612.3152588237999
$ws.Range("S10").Value = 0.6649910641710071
$ws.Range("T10").Value = 0.6649910641710071

# Row 11: update recalculated TPM-derived values
$ws.Range("G11").Value = 16.8273
$ws.Range("H11").Value = 50.4819
$ws.Range("I11").Value = 0.8846442749337277
$ws.Range("J11").Value = 0.8846442749337278
$ws.Range("M11").Value = 0.4282866666666667
$ws.Range("N11").Value = 1.28486
$ws.Range("O11").Value = 0.07962758736516451
$ws.Range("P11").Value = 0.07962758736516451
$ws.Range("Q11").Value = 7.206908225999999
$ws.Range("R11").Value = 64.86217403400001
$ws.Range("S11").Value = 0.07044208928937801
$ws.Range("T11").Value = 0.07044208928937802

# Row 12: update recalculated TPM-derived values
$ws.Range("G12").Value = 16.8273
$ws.Range("H12").Value = 50.4819
$ws.Range("I12").Value = 0.8846442749337277
$ws.Range("J12").Value = 0.8846442749337278
$ws.Range("M12").Value = 0.3522683333333333
$ws.Range("N12").Value = 1.056805
$ws.Range("O12").Value = 0.06549416470700518
$ws.Range("P12").Value = 0.06549416470700517
$ws.Range("Q12").Value = 5.9277249255
$ws.Range("R12").Value = 53.34952432949999
$ws.Range("S12").Value = 0.05793903784961874
$ws.Range("T12").Value = 0.05793903784961874

# Row 13: update recalculated TPM-derived values
$ws.Range("G13").Value = 16.8273
$ws.Range("H13").Value = 50.4819
$ws.Range("I13").Value = 0.8846442749337277
$ws.Range("J13").Value = 0.8846442749337278
$ws.Range("M13").Value = 0.5549326666666667
$ws.Range("N13").Value = 1.664798
$ws.Range("O13").Value = 0.1031737684964519
$ws.Range("P13").Value = 0.1031737684964519
$ws.Range("Q13").Value = 9.338018461799999
$ws.Range("R13").Value = 84.04216615620001
$ws.Range("S13").Value = 0.09127208362372395
$ws.Range("T13").Value = 0.09127208362372394

# Row 14: new Resolving-Mac sending-cluster row
$ws.Range("A14").Value = "Resolving-Mac"
$ws.Range("B14").Value = "Sema3a"
$ws.Range("C14").Value = "Plxna4"
$ws.Range("D14").Value = "ECs"
$ws.Range("E14").Value = 1
$ws.Range("F14").Value = 0.3333333333333333
$ws.Range("G14").Value = 0.016138
$ws.Range("H14").Value = 0.048414
$ws.Range("I14").Value = 0.0008484064174811467
$ws.Range("J14").Value = 0.0008484064174811467
$ws.Range("K14").Value = 3
$ws.Range("L14").Value = 1
$ws.Range("M14").Value = 4.043133999999999
$ws.Range("N14").Value = 12.129402
$ws.Range("O14").Value = 0.7517044794313785
$ws.Range("P14").Value = 0.7517044794313784
$ws.Range("Q14").Value = 0.06524809649199999
$ws.Range("R14").Value = 0.5872328684279999
$ws.Range("S14").Value = 0.0006377509043989061
$ws.Range("T14").Value = 0.0006377509043989061

# Row 15: new Resolving-Mac sending-cluster row
$ws.Range("A15").Value = "Resolving-Mac"
$ws.Range("B15").Value = "Sema3a"
$ws.Range("C15").Value = "Plxna4"
$ws.Range("D15").Value = "FAPs"
$ws.Range("E15").Value = 1
$ws.Range("F15").Value = 0.3333333333333333
$ws.Range("G15").Value = 0.016138
$ws.Range("H15").Value = 0.048414
$ws.Range("I15").Value = 0.0008484064174811467
$ws.Range("J15").Value = 0.0008484064174811467
$ws.Range("K15").Value = 3
$ws.Range("L15").Value = 1
$ws.Range("M15").Value = 0.4282866666666667
$ws.Range("N15").Value = 1.28486
$ws.Range("O15").Value = 0.07962758736516451
$ws.Range("P15").Value = 0.07962758736516451
$ws.Range("Q15").Value = 0.006911690226666667
$ws.Range("R15").Value = 0.06220521204000001
$ws.Range("S15").Value = 0.00006755655612914624
$ws.Range("T15").Value = 0.00006755655612914624

# Row 16: new Resolving-Mac sending-cluster row
$ws.Range("A16").Value = "Resolving-Mac"
$ws.Range("B16").Value = "Sema3a"
$ws.Range("C16").Value = "Plxna4"
$ws.Range("D16").Value = "MuSCs"
$ws.Range("E16").Value = 1
$ws.Range("F16").Value = 0.3333333333333333
$ws.Range("G16").Value = 0.016138
$ws.Range("H16").Value = 0.048414
$ws.Range("I16").Value = 0.0008484064174811467
$ws.Range("J16").Value = 0.0008484064174811467
$ws.Range("K16").Value = 3
$ws.Range("L16").Value = 1
$ws.Range("M16").Value = 0.3522683333333333
$ws.Range("N16").Value = 1.056805
$ws.Range("O16").Value = 0.06549416470700518
$ws.Range("P16").Value = 0.06549416470700517
$ws.Range("Q16").Value = 0.005684906363333333
$ws.Range("R16").Value = 0.05116415727
$ws.Range("S16").Value = 0.00005556566964499042
$ws.Range("T16").Value = 0.00005556566964499041

# Row 17: new Resolving-Mac sending-cluster row
$ws.Range("A17").Value = "Resolving-Mac"
$ws.Range("B17").Value = "Sema3a"
$ws.Range("C17").Value = "Plxna4"
$ws.Range("D17").Value = "Resolving-Mac"
$ws.Range("E17").Value = 1
$ws.Range("F17").Value = 0.3333333333333333
$ws.Range("G17").Value = 0.016138
$ws.Range("H17").Value = 0.048414
$ws.Range("I17").Value = 0.0008484064174811467
$ws.Range("J17").Value = 0.0008484064174811467
$ws.Range("K17").Value = 3
$ws.Range("L17").Value = 1
$ws.Range("M17").Value = 0.5549326666666667
$ws.Range("N17").Value = 1.664798
$ws.Range("O17").Value = 0.1031737684964519
$ws.Range("P17").Value = 0.1031737684964519
$ws.Range("Q17").Value = 0.008955503374666668
$ws.Range("R17").Value = 0.080599530372
$ws.Range("S17").Value = 0.00008753328730810393
$ws.Range("T17").Value = 0.00008753328730810392
